# Weekly update: shift existing Chirimoya price rows down by inserting a
# new week's worth of data (3 rows) at the top of the data block (row 7),
# then populate those new rows with the latest week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at rows 7-9 (pushes the existing rows 7-29 down to 10-32,
# inheriting the row-above formatting, including the date NumberFormat on column D).
$ws.Range("A7:A9").EntireRow.Insert()

# Row 7: Especial
$ws.Cells.Item(7, 1).Value = 2
$ws.Cells.Item(7, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(7, 3).Value = "Coquimbo"
$ws.Cells.Item(7, 4).Value = 44462
$ws.Cells.Item(7, 5).Value = 4
$ws.Cells.Item(7, 6).Value = "Fruta"
$ws.Cells.Item(7, 7).Value = 100107
$ws.Cells.Item(7, 8).Value = "Otros"
$ws.Cells.Item(7, 9).Value = 100107002
$ws.Cells.Item(7, 10).Value = "Chirimoya"
$ws.Cells.Item(7, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(7, 12).Value = "Especial"
$ws.Cells.Item(7, 13).Value = 300
$ws.Cells.Item(7, 14).Value = 2400
$ws.Cells.Item(7, 15).Value = 2500
$ws.Cells.Item(7, 16).Value = 2450
$ws.Cells.Item(7, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(7, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(7, 19).Value = 2450
$ws.Cells.Item(7, 20).Value = 1

# Row 8: Primera
$ws.Cells.Item(8, 1).Value = 2
$ws.Cells.Item(8, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(8, 3).Value = "Coquimbo"
$ws.Cells.Item(8, 4).Value = 44462
$ws.Cells.Item(8, 5).Value = 4
$ws.Cells.Item(8, 6).Value = "Fruta"
$ws.Cells.Item(8, 7).Value = 100107
$ws.Cells.Item(8, 8).Value = "Otros"
$ws.Cells.Item(8, 9).Value = 100107002
$ws.Cells.Item(8, 10).Value = "Chirimoya"
$ws.Cells.Item(8, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(8, 12).Value = "Primera"
$ws.Cells.Item(8, 13).Value = 400
$ws.Cells.Item(8, 14).Value = 2100
$ws.Cells.Item(8, 15).Value = 2200
$ws.Cells.Item(8, 16).Value = 2150
$ws.Cells.Item(8, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(8, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(8, 19).Value = 2150
$ws.Cells.Item(8, 20).Value = 1

# Row 9: Segunda
$ws.Cells.Item(9, 1).Value = 2
$ws.Cells.Item(9, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(9, 3).Value = "Coquimbo"
$ws.Cells.Item(9, 4).Value = 44462
$ws.Cells.Item(9, 5).Value = 4
$ws.Cells.Item(9, 6).Value = "Fruta"
$ws.Cells.Item(9, 7).Value = 100107
$ws.Cells.Item(9, 8).Value = "Otros"
$ws.Cells.Item(9, 9).Value = 100107002
$ws.Cells.Item(9, 10).Value = "Chirimoya"
$ws.Cells.Item(9, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(9, 12).Value = "Segunda"
$ws.Cells.Item(9, 13).Value = 300
$ws.Cells.Item(9, 14).Value = 1800
$ws.Cells.Item(9, 15).Value = 1900
$ws.Cells.Item(9, 16).Value = 1850
$ws.Cells.Item(9, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(9, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(9, 19).Value = 1850
$ws.Cells.Item(9, 20).Value = 1
